# Update the NBA roster sheet (yahoo/NSY.xlsx) to the post-trade-deadline
# player/position/team listing. Player names stay in column A in their
# existing row order, but two players were swapped out (Toumani Camara,
# who left the roster shown, and Portland Trail Blazers with him) for one
# new player (Matas Buzelis / Chicago Bulls), and every row's Position
# (col B) / Team (col C) pair was recomputed so it correctly matches that
# row's player.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jordan Poole",            "PG,SG", "Washington Wizards"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Jamal Murray",            "PG,SG", "Denver Nuggets"),
    @("Kyrie Irving",            "PG,SG", "Dallas Mavericks"),
    @("Tobias Harris",           "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen",         "SF,PF", "Utah Jazz"),
    @("Zach LaVine",             "SG,SF", "Sacramento Kings"),
    @("Jordan Clarkson",         "SG,SF", "Utah Jazz"),
    @("Bam Adebayo",             "PF,C",  "Miami Heat"),
    @("John Collins",            "PF,C",  "Utah Jazz"),
    @("Kyle Kuzma",              "PF",    "Milwaukee Bucks"),
    @("Gradey Dick",             "SG,SF", "Toronto Raptors"),
    @("Joel Embiid",             "C",     "Philadelphia 76ers"),
    @("Matas Buzelis",           "SF,PF", "Chicago Bulls"),
    @("Rui Hachimura",           "SF,PF", "Los Angeles Lakers"),
    @("CJ McCollum",             "PG,SG", "New Orleans Pelicans")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
